$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new recipient row (row 3): same email as row 2, new name + subject
$ws.Range("A3").Value = "fc.krkim@gmail.com"
$ws.Range("B3").Value = "김영환"
$ws.Range("C3").Value = "[패스트몰] 2022-09-22 상품발주 확인 요청9"

# Turn A3 into a mailto hyperlink, like A2
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:fc.krkim@gmail.com") | Out-Null

# Hyperlinks.Add re-styles the cell with its own copy of the hyperlink style;
# put it back on the same named "하이퍼링크" style A2 already uses.
$ws.Range("A3").Style = $ws.Range("A2").Style

# Update the sheet view: drop the B1 frozen/topLeft offset, zoom out to 175%,
# and move the active selection to D3
$ws.Range("D3").Select() | Out-Null
$excel.ActiveWindow.Zoom = 175
